# fixed | i-475 | Se agrega campos nombre - nombre secundario importacion productos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column A header from "Descripción" to "Nombre"
$ws.Range("A1").Value = "Nombre"

# Add new column O: "Descripcion"
$ws.Range("O1").Value = "Descripcion"

# Add new column P: "Nombre secundario"
$ws.Range("P1").Value = "Nombre secundario"

# Fill in column O values
$ws.Range("O2").Value = "desc 1"
$ws.Range("O3").Value = "desc 2"

# Fill in column P values
$ws.Range("P2").Value = "nombre sec 11"
$ws.Range("P3").Value = "nombre sec 22"

# Best-fit column P to its widest entry (mirrors Excel's auto column sizing
# after the new data was typed in)
$ws.Columns.Item(16).ColumnWidth = 17.71

# Update selection cursor to O8, matching the target workbook state
$ws.Range("O8").Select()
